# Update the part number for row 4 (Bourns potentiometer) from the
# 4015U variant to the 4025U variant.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cart")

$ws.Range("A4").Value = "652-PTV09A-4025UB103"
$ws.Range("B4").Value = "PTV09A-4025U-B103"
$ws.Range("G4").Value = "652-PTV09A-4025UB103"

# Update the sheet view selection/scroll position to match the saved state.
$ws.Range("F21").Select()
$ws.Application.ActiveWindow.ScrollColumn = 4
